$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3,5,7,9,11,13,15,17 had their G and H column values removed (cleared)
foreach ($r in 3,5,7,9,11,13,15,17) {
    $ws.Range("G$r`:H$r").ClearContents()
}
